# Fixes during Regression Testing
# Update the "DateProd" (col B) notes timestamps and, for rows whose
# prod run failed, flip "ResultDemo" (col C) from Pass to Fail.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [string]$DateProd,
        [string]$ResultDemo
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 2).Value = $DateProd
    if ($ResultDemo) {
        $ws.Cells.Item($Row, 3).Value = $ResultDemo
    }
}

# VT-P-DebitVoid-DualCF-Generic
Set-Row "VT-P-DebitVoid-DualCF-Generic" 2 "Thu Nov 20 20:29:07 IST 2025" $null

# VT-P-DebitVoid-SingleCF-Generic
Set-Row "VT-P-DebitVoid-SingleCF-Generic" 2 "Thu Nov 20 20:35:23 IST 2025" $null
Set-Row "VT-P-DebitVoid-SingleCF-Generic" 3 "Thu Nov 20 20:36:48 IST 2025" $null
Set-Row "VT-P-DebitVoid-SingleCF-Generic" 4 "Thu Nov 20 20:37:51 IST 2025" $null
Set-Row "VT-P-DebitVoid-SingleCF-Generic" 5 "Thu Nov 20 20:39:20 IST 2025" $null

# VT-P-DebitVoid-NoCF-Generic
Set-Row "VT-P-DebitVoid-NoCF-Generic" 2 "Thu Nov 20 20:30:44 IST 2025" $null
Set-Row "VT-P-DebitVoid-NoCF-Generic" 3 "Thu Nov 20 20:31:56 IST 2025" $null
Set-Row "VT-P-DebitVoid-NoCF-Generic" 4 "Thu Nov 20 20:33:05 IST 2025" $null
Set-Row "VT-P-DebitVoid-NoCF-Generic" 5 "Thu Nov 20 20:34:12 IST 2025" $null

# VT-P-DebitCredit-DualCF-Generic
Set-Row "VT-P-DebitCredit-DualCF-Generic" 2 "Thu Nov 20 20:01:48 IST 2025" $null
Set-Row "VT-P-DebitCredit-DualCF-Generic" 3 "Thu Nov 20 20:02:57 IST 2025" $null
Set-Row "VT-P-DebitCredit-DualCF-Generic" 4 "Thu Nov 20 20:04:12 IST 2025" $null
Set-Row "VT-P-DebitCredit-DualCF-Generic" 5 "Thu Nov 20 20:05:19 IST 2025" $null

# VT-P-DebitCredit-SingleCF-Gener (rows 2-4 flip ResultDemo to Fail)
Set-Row "VT-P-DebitCredit-SingleCF-Gener" 2 "Thu Nov 20 20:11:01 IST 2025" "Fail"
Set-Row "VT-P-DebitCredit-SingleCF-Gener" 3 "Thu Nov 20 20:12:17 IST 2025" "Fail"
Set-Row "VT-P-DebitCredit-SingleCF-Gener" 4 "Thu Nov 20 20:13:33 IST 2025" "Fail"
Set-Row "VT-P-DebitCredit-SingleCF-Gener" 5 "Thu Nov 20 16:13:01 IST 2025" $null

# VT-P-DebitCredit-NoCF-Generic
Set-Row "VT-P-DebitCredit-NoCF-Generic" 2 "Thu Nov 20 20:06:26 IST 2025" $null
Set-Row "VT-P-DebitCredit-NoCF-Generic" 3 "Thu Nov 20 20:07:30 IST 2025" $null
Set-Row "VT-P-DebitCredit-NoCF-Generic" 4 "Thu Nov 20 20:08:38 IST 2025" $null
Set-Row "VT-P-DebitCredit-NoCF-Generic" 5 "Thu Nov 20 20:09:51 IST 2025" $null

# VT-C-DebitCredit-DualCF-Generic
Set-Row "VT-C-DebitCredit-DualCF-Generic" 2 "Thu Nov 20 19:57:47 IST 2025" $null

# VT-C-DebitCredit-SingleCF-Gener (row 2 flips ResultDemo to Fail)
Set-Row "VT-C-DebitCredit-SingleCF-Gener" 2 "Thu Nov 20 22:32:39 IST 2025" "Fail"

# VT-C-DebitCredit-NoCF-Generic
Set-Row "VT-C-DebitCredit-NoCF-Generic" 2 "Thu Nov 20 19:59:07 IST 2025" $null
